# The "请假" (leave-request) roster sheet is being reset to a blank
# template: drop the old class-roster rows (names / student IDs / "lab1"
# header) and replace them with a fresh 3-column header row
# (姓名 / 学号 / 请假日期).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old A1:B1 title cell was merged ("lab1") - unmerge before touching
# the cells underneath it.
$ws.Range("A1:B1").UnMerge()

# Wipe every old value (names, student ids, the "lab1"/"?" header row)
# and its formatting so nothing is left behind.
$ws.Cells.Clear()

# New header row.
$ws.Range("A1").Value = "姓名"
$ws.Range("B1").Value = "学号"
$ws.Range("C1").Value = "请假日期"

# Leave the cursor where the saved workbook shows it.
$ws.Range("E6").Select()
